$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A3 value (9790581357 -> 13456)
$ws.Range("A3").Value = 13456

# Add new row with A4 = 999
$ws.Range("A4").Value = 999

# Bump the sheet's max row outline level to 3 (outlineLevelRow 2 -> 3) without
# leaving a visible outlineLevel on any surviving row: group a scratch row
# beyond the used range at level 3, then remove it again.
$ws.Rows("5:5").OutlineLevel = 3
$ws.Rows("5:5").Delete()

# Move the active selection to the new last cell A4
$ws.Range("A4").Select()
